$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`nLipids and endothelium are pivotal players on the scene of atherosclerosis and their interaction is crucial for the establishment of the pathological processes.`n The endothelium is not only the border of the arterial wall: it plays a key role in regulating circulating fatty acids and lipoproteins and vice versa it is regulated by these lipidic molecules thereby promoting atherosclerosis.`n Inflammation is another important element in the relationship between lipids and endothelium.`n Recently, proprotein convertase subtilisin/kexin type 9 (PCSK9) has been recognized as a fundamental regulator of LDL-C and anti-PCSK9 monoclonal antibodies have been approved for therapeutic use in hypercholesterolemia, with the promise to subvert the natural history of the disease.`n Moreover, growing experimental and clinical evidence is enlarging our understanding of the mechanisms through which this protein may facilitate the genesis of atherosclerosis, independently of its impact on lipid metabolism.`n In addition, environmental stimuli may affect the post-transcriptional regulation of genes through micro-RNAs, which in turn play a key role in orchestrating the crosstalk between endothelium and cholesterol.`n Advances in experimental research, with development of high throughput techniques, have led, over the last century, to a tremendous progress in the understanding and fine tuning of the molecular mechanisms leading to atherosclerosis.`n Identification of pivotal keystone molecules bridging lipid metabolism, endothelial dysfunction and atherogenesis will provide the mechanistic substrate to test valuable targets for prediction, prevention and treatment of atherosclerosis-related disease.`n"
$ws.Range("E2").Value = "[ D.%D’Ardes%null%0,  F.%Santilli%null%1,  M. T.%Guagnano%null%1,  M.%Bucci%null%1,  F.%Cipollone%null%1]"
$ws.Range("E3").Value = "[Chen%Chen%NULL%0,           Feng%Wang%NULL%0,           Peng%Chen%NULL%1,           Jiangang%Jiang%NULL%0,           Guanglin%Cui%NULL%1,           Ning%Zhou%NULL%0,           Francesco%Moroni%NULL%1,           Javid J.%Moslehi%NULL%0,           Enrico%Ammirati%enrico.ammirati@ospedaleniguarda.it%0,           Dao Wen%Wang%dwwang@tjh.tjmu.edu.cn%0]"
$ws.Range("E4").Value = "[Marcello%Covino%NULL%0,           Giuseppe%De Matteis%giuseppe.dematteis@policlinicogemelli.it%1,           Maria Livia%Burzo%NULL%2,           Maria Livia%Burzo%NULL%0,           Michele%Santoro%NULL%1,           Mariella%Fuorlo%NULL%1,           Luca%Sabia%NULL%1,           Claudio%Sandroni%NULL%1,           Antonio%Gasbarrini%NULL%0,           Francesco%Franceschi%NULL%1,           Giovanni%Gambassi%NULL%1,           NULL%NULL%NULL%0]"
$ws.Range("E5").Value = "[C%Felice%NULL%1,           C%Nardin%NULL%1,           G L%Di Tanna%NULL%1,           U%Grossi%NULL%1,           E%Bernardi%NULL%1,           L%Scaldaferri%NULL%1,           M%Romagnoli%NULL%1,           L%Tonon%NULL%1,           P%Cavasin%NULL%1,           S%Novello%NULL%1,           R%Scarpa%NULL%1,           A%Farnia%NULL%1,           E%De Menis%NULL%1,           R%Rigoli%NULL%1,           F%Cinetto%NULL%1,           P%Pauletto%NULL%1,           C%Agostini%NULL%1,           M%Rattazzi%marcello.rattazzi@unipd.it%1]"
$ws.Range("E6").Value = "[Zheyong%Huang%NULL%0,           Jiatian%Cao%NULL%1,           Yumeng%Yao%NULL%1,           Xuejuan%Jin%NULL%1,           Zhe%Luo%NULL%1,           Yuan%Xue%NULL%1,           Chouwen%Zhu%NULL%1,           Yanan%Song%NULL%1,           Ying%Wang%NULL%0,           Yunzeng%Zou%NULL%1,           Juying%Qian%NULL%1,           Kaihuan%Yu%NULL%1,           Hui%Gong%NULL%1,           Junbo%Ge%NULL%1]"
$ws.Range("E7").Value = "[Katherine W%Lam%NULL%0,           Kenneth W%Chow%NULL%1,           Jonathan%Vo%NULL%1,           Wei%Hou%NULL%0,           Haifang%Li%NULL%0,           Paul S%Richman%NULL%1,           Sandeep K%Mallipattu%NULL%1,           Hal A%Skopicki%NULL%1,           Adam J%Singer%NULL%1,           Tim Q%Duong%tim.duong@stonybrook.edu%1]"
$ws.Range("E8").Value = "[Xiulan%Liu%NULL%0,           Yi%Liu%NULL%0,           Keliang%Chen%NULL%1,           Suying%Yan%NULL%1,           Xiangrong%Bai%NULL%1,           Juan%Li%947281063@qq.com%0,           Dong%Liu%ld2069@outlook.com%0,           Dong%Liu%ld2069@outlook.com%0]"
$ws.Range("D9").Value = "`nSince the beginning of the coronavirus disease 2019 (COVID-19) outbreak initiated on the Diamond Princess Cruise Ship at Yokohama harbor in February 2020, we have been doing our best to treat COVID-19 patients.`n In animal experiments, angiotensin converting enzyme inhibitors (ACEIs) and angiotensin II type-1 receptor blockers (ARBs) are reported to suppress the downregulation of angiotensin converting enzyme 2 (ACE2), and they may inhibit the worsening of pathological conditions.`n We aimed to examine whether preceding use of ACEIs and ARBs affected the clinical manifestations and prognosis of COVID-19 patients.`n One hundred fifty-one consecutive patients (mean age 60 ± 19 years) with polymerase-chain-reaction proven severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) infection who were admitted to six hospitals in Kanagawa Prefecture, Japan, were analyzed in this multicenter retrospective observational study.`n Among all COVID-19 patients, in the multiple regression analysis, older age (age ≥ 65 years) was significantly associated with the primary composite outcome (odds ratio (OR) 6.63, 95% confidence interval (CI) 2.28–22.78, P < 0.001), which consisted of (i) in-hospital death, (ii) extracorporeal membrane oxygenation, (iii) mechanical ventilation, including invasive and noninvasive methods, and (iv) admission to the intensive care unit.`n In COVID-19 patients with hypertension, preceding ACEI/ARB use was significantly associated with a lower occurrence of new-onset or worsening mental confusion (OR 0.06, 95% CI 0.002–0.69, P = 0.02), which was defined by the confusion criterion, which included mild disorientation or hallucination with an estimation of medical history of mental status, after adjustment for age, sex, and diabetes.`n In conclusion, older age was a significant contributor to a worse prognosis in COVID-19 patients, and ACEIs/ARBs could be beneficial for the prevention of confusion in COVID-19 patients with hypertension.`n"
$ws.Range("E9").Value = "[ Yasushi%Matsuzawa%null%0,  Hisao%Ogawa%null%0,  Kazuo%Kimura%null%0,  Masaaki%Konishi%null%0,  Jin%Kirigaya%null%0,  Kazuki%Fukui%null%0,  Kengo%Tsukahara%null%0,  Hiroyuki%Shimizu%null%0,  Keisuke%Iwabuchi%null%0,  Yu%Yamada%null%0,  Kenichiro%Saka%null%0,  Ichiro%Takeuchi%null%0,  Toshio%Hirano%null%0,  Kouichi%Tamura%null%0]"
$ws.Range("E10").Value = "[Juan%Meng%NULL%0,           Guohui%Xiao%NULL%1,           Juanjuan%Zhang%NULL%1,           Xing%He%NULL%1,           Min%Ou%NULL%1,           Jing%Bi%NULL%1,           Rongqing%Yang%NULL%1,           Wencheng%Di%NULL%1,           Zhaoqin%Wang%NULL%1,           Zigang%Li%NULL%1,           Hong%Gao%NULL%0,           Lei%Liu%NULL%0,           Guoliang%Zhang%NULL%1]"
$ws.Range("D11").Value = "`nUltrasound stimulation is an emerging noninvasive option in treating neuropsychiatric disorders.`n The present study investigates the behavioral alterations resulting from ultrasound stimulation on the nucleus accumbens (NAc) in freely moving mice.`n Our results show that an acute ultrasound stimulation on the NAc, rather than the visual cortex or auditory cortex, led to a pronounced avoidance behavior, while repeated NAc ultrasound stimulation resulted in an obvious conditioned place aversion with changes in synaptic protein (GluA1/2 subunit) expression.`n Notably, NAc ultrasound stimulation suppressed the morphine-induced conditioned place preference.`n The results provide evidence that NAc ultrasound stimulation can be applied as a potential noninvasive therapeutic option in treating psychiatric disorders.`n"
$ws.Range("E11").Value = "[ Lili%Niu%null%0,  Yanchen%Guo%null%1,  Zhengrong%Lin%null%1,  Zhe%Shi%null%1,  Tianyuan%Bian%null%1,  Lin%Qi%null%1,  Long%Meng%null%1,  Anthony A.%Grace%null%1,  Hairong%Zheng%null%1,  Ti-Fei%Yuan%null%1]"
$ws.Range("E12").Value = "[Wei%Pan%xref no email%0,    Jishou%Zhang%xref no email%0,    Menglong%Wang%xref no email%0,    Jing%Ye%xref no email%0,    Yao%Xu%xref no email%0,    Bo%Shen%xref no email%0,    Hua%He%xref no email%0,    Zhen%Wang%xref no email%0,    Di%Ye%xref no email%0,    Mengmeng%Zhao%xref no email%0,    Zhen%Luo%xref no email%0,    Mingxiao%Liu%xref no email%0,    Pingan%Zhang%xref no email%0,    Jian%Gu%xref no email%0,    Menglin%Liu%xref no email%0,    Dan%Li%xref no email%0,    Jianfang%Liu%xref no email%0,    Jun%Wan%xref no email%0]"
$ws.Range("E13").Value = "[Celestino%Sardu%drsarducele@gmail.com%0,           Paolo%Maggi%NULL%2,           Paolo%Maggi%NULL%0,           Vincenzo%Messina%NULL%1,           Pasquale%Iuliano%NULL%2,           Pasquale%Iuliano%NULL%0,           Antonio%Sardu%NULL%1,           Vincenzo%Iovinella%NULL%1,           Giuseppe%Paolisso%NULL%1,           Raffaele%Marfella%NULL%1]"
$ws.Range("E14").Value = "[International%Society of Hypertension%coreGivesNoEmail%1,          Li%J%coreGivesNoEmail%1,          Zhang%P%coreGivesNoEmail%1]"
$ws.Range("E15").Value = "[Naci%Şenkal%NULL%0,           Rasimcan%Meral%NULL%1,           Alpay%Medetalibeyoğlu%NULL%1,           Hilal%Konyaoğlu%NULL%1,           Murat%Köse%NULL%1,           Tufan%Tükek%NULL%1]"
$ws.Range("E16").Value = "[Minkook%Son%xref no email%1,    Jeongkuk%Seo%xref no email%1,    Sung%Yang%xref no email%1]"
$ws.Range("E17").Value = "[Zhongchao%Wang%NULL%1,           Dewei%Zhang%NULL%1,           Shengming%Wang%NULL%1,           Yanhua%Jin%NULL%1,           Jianbo%Huan%NULL%1,           Yue%Wu%NULL%0,           Cheng%Xia%NULL%1,           Zhe%Li%NULL%1,           Xingshun%Qi%NULL%1,           Duanzhen%Zhang%NULL%1,           Xiumin%Han%NULL%1,           Xianyang%Zhu%NULL%1,           Ying%Qu%NULL%1,           Qiguang%Wang%NULL%1]"
$ws.Range("E18").Value = "[Honglu%Zhou%NULL%1,           Songmei%Wang%NULL%1,           Lorenz%von Seidlein%NULL%1,           Xuanyi%Wang%xywang@shmu.edu.cn%1]"
$ws.Range("E19").Value = "[Guang%Yang%xref no email%0,    Zihu%Tan%xref no email%0,    Ling%Zhou%xref no email%0,    Min%Yang%xref no email%0,    Lang%Peng%xref no email%0,    Jinjin%Liu%xref no email%0,    Jingling%Cai%xref no email%0,    Ru%Yang%xref no email%0,    Junyan%Han%xref no email%0,    Yafei%Huang%xref no email%0,    Shaobin%He%xref no email%0]"
$ws.Range("E20").Value = "[Yuan%Yuan%NULL%0,           Dan%Liu%NULL%0,           Shaoqing%Zeng%NULL%1,           Siyuan%Wang%NULL%1,           Sen%Xu%NULL%0,           Ya%Wang%NULL%0,           Ruidi%Yu%NULL%1,           Yue%Gao%NULL%0,           Huayi%Li%NULL%1,           Xinxia%Feng%NULL%1,           Ning%Zhou%NULL%0,           Chunxia%Zhao%NULL%1,           Qinglei%Gao%NULL%1]"
$ws.Range("E21").Value = "[Peng%Zhang%NULL%0,           Lihua%Zhu%NULL%0,           Jingjing%Cai%NULL%0,           Fang%Lei%NULL%0,           Juan-Juan%Qin%NULL%0,           Jing%Xie%NULL%2,           Ye-Mao%Liu%NULL%0,           Yan-Ci%Zhao%NULL%1,           Xuewei%Huang%NULL%0,           Lijin%Lin%NULL%0,           Meng%Xia%NULL%0,           Ming-Ming%Chen%NULL%0,           Xu%Cheng%NULL%1,           Xiao%Zhang%NULL%0,           Deliang%Guo%NULL%1,           Yuanyuan%Peng%NULL%1,           Yan-Xiao%Ji%NULL%0,           Jing%Chen%NULL%0,           Zhi-Gang%She%NULL%0,           Yibin%Wang%NULL%0,           Qingbo%Xu%NULL%0,           Renfu%Tan%NULL%1,           Haitao%Wang%NULL%1,           Jun%Lin%NULL%2,           Pengcheng%Luo%NULL%1,           Shouzhi%Fu%NULL%0,           Hongbin%Cai%NULL%1,           Ping%Ye%NULL%0,           Bing%Xiao%NULL%1,           Weiming%Mao%NULL%0,           Liming%Liu%NULL%1,           Youqin%Yan%NULL%2,           Mingyu%Liu%NULL%1,           Manhua%Chen%NULL%1,           Xiao-Jing%Zhang%NULL%0,           Xinghuan%Wang%NULL%0,           Rhian M.%Touyz%NULL%1,           Jiahong%Xia%NULL%0,           Bing-Hong%Zhang%NULL%0,           Xiaodong%Huang%NULL%0,           Yufeng%Yuan%NULL%0,           Rohit%Loomba%NULL%1,           Peter P.%Liu%NULL%1,           Hongliang%Li%NULL%0]"
$ws.Range("E22").Value = "[Xian%Zhou%NULL%1,           Jingkang%Zhu%NULL%1,           Tao%Xu%NULL%0]"
